$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 6 new rows (49-54) below the current last row (48). ---
# Insert() inherits the format of the row above (48), which is already the
# "4/5" (non-separator) style, so most of the new rows get the right style
# for free; rows 50/52 (continuation lines) and the now-not-last row 48
# (separator look) are fixed up below via format-only paste.
$ws.Rows("49:54").Insert()

# --- Re-point row 48's style: it used to be the last (non-separator) row,
# now row 48 is a mid-table row that should look like the other "divider"
# rows (e.g. 43/46/47), i.e. cellXf 8/9 instead of 4/5. ---
$ws.Range("A43:E43").Copy()
$ws.Range("A48:E48").PasteSpecial(-4122)

# --- Rows 50 and 52 are "continuation" rows (2nd line of a 2-line entry),
# which use cellXf 6/7 (matches rows such as 42/45). ---
$ws.Range("A42:E42").Copy()
$ws.Range("A50:E50").PasteSpecial(-4122)
$ws.Range("A42:E42").Copy()
$ws.Range("A52:E52").PasteSpecial(-4122)

# --- Row 54's first column has no cell at all (same pattern as rows 5/6/9/10). ---
$ws.Range("A54").Clear()

# --- New cell values (filenames, row numbers, EN/RU/garbled-RU text). ---
# NB: throughout this sheet every apostrophe in the source text is stored
# as a literal backslash followed by an apostrophe (two characters, not an
# escape sequence) -- e.g. "It\'s" -- and "\n" is likewise a literal
# backslash+n, not a real newline. The doubled '' below is plain
# PowerShell single-quote escaping for the literal backslash-apostrophe.
$ws.Range("A49").Value = 'SCRIPT/G01P03A/um2101.ssb'
$ws.Range("B49").Value = 357
$ws.Range("C49").Value = ' I\''m happy to see your team back,\n[hero]!'
$ws.Range("D49").Value = ' Я рада, что твоя команда снова\nв строю, [hero]!'
$ws.Range("E49").Value = ' Ÿ ñàäà, œóï óâïÿ ëïíàîäà òîïâà\nâ òóñïý, [hero]!'

$ws.Range("B50").Value = 360
$ws.Range("C50").Value = ' I hope we can work together!'
$ws.Range("D50").Value = ' Надеюсь, мы ещё будем работать\nвместе!'
$ws.Range("E50").Value = ' Îàäåýòû, íú åþæ áôäåí ñàáïóàóû\nâíåòóå!'

$ws.Range("A51").Value = 'SCRIPT/G01P03A/um2201.ssb'
$ws.Range("B51").Value = 329
$ws.Range("C51").Value = ' I\''m going out to [CS:P]Brine Cave[CR]\ntomorrow too!'
$ws.Range("D51").Value = ' Завтра я тоже отправлюсь в\n[CS:P]Пещеру у Моря[CR]!'
$ws.Range("E51").Value = ' Èàâóñà ÿ óïçå ïóðñàâìýòû â\n[CS:P]Ðåþåñô ô Íïñÿ[CR]!'

$ws.Range("B52").Value = 332
$ws.Range("C52").Value = ' Let\''s do good!'
$ws.Range("D52").Value = ' Мы справимся!'
$ws.Range("E52").Value = ' Íú òðñàâéíòÿ!'

$ws.Range("A53").Value = 'SCRIPT/G01P03A/um2402.ssb'
$ws.Range("B53").Value = 301
$ws.Range("C53").Value = ' It\''s not possible to add team\nmembers when exploring the [CS:P]Hidden Land[CR].[K]\nThat\''s very disappointing…'
$ws.Range("D53").Value = ' Во время исследования [CS:P]Сокрытых\nЗемель[CR] нельзя брать с собой других членов\nгруппы.[K] Это расстраивает...'
$ws.Range("E53").Value = ' Âï âñåíÿ éòòìåäïâàîéÿ [CS:P]Òïëñúóúö\nÈåíåìû[CR] îåìûèÿ áñàóû ò òïáïê äñôãéö œìåîïâ\nãñôððú.[K] Üóï ñàòòóñàéâàåó...'

$ws.Range("B54").Value = 304
$ws.Range("C54").Value = ' Don\''t give up, though! You have\nto succeed!'
$ws.Range("D54").Value = ' Но не сдавайтесь! У вас всё\nполучится!'
$ws.Range("E54").Value = ' Îï îå òäàâàêóåòû! Ô âàò âòæ\nðïìôœéóòÿ!'

# --- Row heights for the new rows (wrapped-text rows sized to their content). ---
$ws.Rows(49).RowHeight = 43.2
$ws.Rows(50).RowHeight = 21.6
$ws.Rows(51).RowHeight = 43.2
$ws.Rows(53).RowHeight = 52.2
$ws.Rows(54).RowHeight = 21.6

# --- Scroll/selection state, matching where the sheet was left after editing. ---
$ws.Range("D54").Select()
